$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.433.81"
$ws.Range("E2").Value = "  +2.25%  "
$ws.Range("D3").Value = "3.180.64"
$ws.Range("E3").Value = "  +4.06%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "573.22"
$ws.Range("E5").Value = "  +2.75%  "
$ws.Range("D6").Value = "151.99"
$ws.Range("E6").Value = "  +6.57%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.179.03"
$ws.Range("E8").Value = "  +4.04%  "
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  +3.93%  "
$ws.Range("E10").Value = "  +5.87%  "
$ws.Range("D11").Value = "6.21"
$ws.Range("E11").Value = "  +2.43%  "
$ws.Range("D12").Value = "0.508"
$ws.Range("E12").Value = "  +6.43%  "
$ws.Range("D13").Value = "0.0000281"
$ws.Range("E13").Value = "  +21.59%  "
$ws.Range("D14").Value = "38.37"
$ws.Range("E14").Value = "  +9.57%  "
$ws.Range("D15").Value = "3.703.17"
$ws.Range("E15").Value = "  +4.53%  "
$ws.Range("D16").Value = "65.468.03"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").Value = "7.24"
$ws.Range("E17").Value = "  +7.69%  "
$ws.Range("D18").Value = "3.186.36"
$ws.Range("E18").Value = "  +4.28%  "
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").Value = "514.94"
$ws.Range("E20").Value = "  +8.20%  "
$ws.Range("D21").Value = "14.98"
$ws.Range("E21").Value = "  +6.70%  "
$ws.Range("E22").Value = "  +12.47%  "
$ws.Range("D23").Value = "0.740"
$ws.Range("E23").Value = "  +9.31%  "
$ws.Range("D24").Value = "7.89"
$ws.Range("E24").Value = "  +4.35%  "
$ws.Range("D25").Value = "85.13"
$ws.Range("E25").Value = "  +4.39%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "9.10"
$ws.Range("E27").Value = "  +14.77%  "
$ws.Range("D28").Value = "2.92"
$ws.Range("E28").Value = "  +4.48%  "
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  +9.22%  "
$ws.Range("D30").Value = "28.18"
$ws.Range("E30").Value = "  +7.48%  "
$ws.Range("D31").Value = "2.81"
$ws.Range("E31").Value = "  +15.56%  "
$ws.Range("D32").Value = "1.23"
$ws.Range("E32").Value = "  +8.30%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").Value = "6.30"
$ws.Range("E34").Value = "  +12.70%  "
$ws.Range("D35").Value = "6.69"
$ws.Range("E35").Value = "  +8.01%  "
$ws.Range("D36").Value = "55.72"
$ws.Range("E36").Value = "  +1.85%  "
$ws.Range("D37").Value = "479.33"
$ws.Range("E37").Value = "  +8.12%  "
$ws.Range("D38").Value = "3.17"
$ws.Range("E38").Value = "  +12.12%  "
$ws.Range("D39").Value = "0.0880"
$ws.Range("E39").Value = "  +9.26%  "
$ws.Range("D40").Value = "0.0423"
$ws.Range("E40").Value = "  +4.42%  "
$ws.Range("D41").Value = "3.147.24"
$ws.Range("E41").Value = "  +6.15%  "
$ws.Range("D42").Value = "8.66"
$ws.Range("E42").Value = "  +5.56%  "
$ws.Range("E43").Value = "  +7.39%  "
$ws.Range("D44").Value = "2.51"
$ws.Range("E44").Value = "  +16.78%  "
$ws.Range("D45").Value = "0.290"
$ws.Range("E45").Value = "  +11.52%  "
$ws.Range("D46").Value = "29.29"
$ws.Range("E46").Value = "  +6.07%  "
$ws.Range("E47").Value = "  +15.14%  "
$ws.Range("E49").Value = "  +3.10%  "
$ws.Range("E50").Value = "  +13.06%  "
$ws.Range("D51").Value = "123.50"
$ws.Range("E51").Value = "  +5.39%  "
